# -----------------------------------------------------------------------
# Workbook edit: rename "W-J 2013" -> "Other W-J", make it the active
# sheet, fix a leftover "2013!!!" rich-text fragment in K2, append two
# rows of data that were missing, and refresh that sheet's layout
# (uniform column width, row heights, page margins, header/footer).
# -----------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "W-J 2017"  (unchanged)
$ws2 = $wb.Worksheets.Item(2)   # "W-J 2013"  -> "Other W-J"

# 1. Rename the second sheet.
$ws2.Name = "Other W-J"

# 2. K2 on sheet 2 carried a stray rich-text run ("-> 2013!!! ...") that
#    duplicated sheet 1's K2 text; replace it with the plain/clean text.
$ws2.Range("K2").Value = @"
→ (12)15-17 PC, 5+ trefl (naturalne, słaba piątka może też pasować na  1BA – patrz dalej)
<br><font color=red>w wersji pro od 15 a nawet 16 PC!!! - słaby trefl wtedy do BA albo 1 karo (o ile jest 4 karo)</font>
<br>→ 12-17 PC, układ 4414  tzw TRÓJKOLORÓWKA
(singiel karo)
<br>→ 12-14 PC, skład bez atutowy  tzw. PRZYGOTOWAWCZE 
(słabe NT) 
<br>→ 18+ PC,  skład dowolny
(silne dowolne)
uwaga: może też być licytowane przy zał że partner ma skład zrównoważony i jak powiem coś innego na co on spasuje (bez punktów) i ucieknie nam końcówka (np. na układzie 5-4 lub 6-4 – układy na dwóch piątkach raczej 
"@

# 3. Sheet 2 was missing the two data rows that sheet 1 already has in
#    row 5 and row 6 - bring them across (values, then formats).
$ws1.Range("A5:M6").Copy()
$ws2.Range("A5:M6").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws1.Range("A5:M6").Copy()
$ws2.Range("A5:M6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 4. Re-wrap / resize the rows that now hold (slightly) different text.
$ws2.Rows.Item(2).RowHeight = 323
$ws2.Rows.Item(3).RowHeight = 170.3
$ws2.Rows.Item(4).RowHeight = 282.8
$ws2.Rows.Item(5).RowHeight = 129.6
$ws2.Rows.Item(6).RowHeight = 129.6

# 5. Give every column on the sheet the same (narrower) width instead of
#    the old A / C:J / K / rest grouping.
$ws2.Range($ws2.Cells.Item(1, 1), $ws2.Cells.Item(1, 1025)).ColumnWidth = 10.508503401360565

# 6. Page margins (points = inches * 72).
$ws2.PageSetup.LeftMargin   = 0.7875 * 72
$ws2.PageSetup.RightMargin  = 0.7875 * 72
$ws2.PageSetup.TopMargin    = 1.05277777777778 * 72
$ws2.PageSetup.BottomMargin = 1.05277777777778 * 72
$ws2.PageSetup.HeaderMargin = 0.7875 * 72
$ws2.PageSetup.FooterMargin = 0.7875 * 72

# 7. Give the sheet a printed header/footer.
$ws2.PageSetup.CenterHeader = '&"Times New Roman,Normalny"&12&A'
$ws2.PageSetup.CenterFooter = '&"Times New Roman,Normalny"&12Strona &P'

# 8. Sheet 1's selection resets to A1, sheet 2 becomes the active/visible
#    tab at 140% zoom with A1 selected.
$ws1.Range("A1").Select() | Out-Null
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 140
$ws2.Range("A1").Select() | Out-Null
